$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "role" column before "first_name" (which currently is column D).
# This shifts first_name -> E, last_name -> F, and inherits the bold header
# style + border from the neighbouring header cell.
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("D1").Value = "role"

# Give the new rightmost header cell (F1, last_name) a left+right border.
$ws.Range("F1").Borders.Item(7).LineStyle = 1
$ws.Range("F1").Borders.Item(10).LineStyle = 1

# Mark two cells below "email" with the Hyperlink cell style (used later for
# CSV-imported e-mail links), but leave them empty for now.
$ws.Hyperlinks.Add($ws.Range("C2"), "http://example.com")
$ws.Range("C2").ClearContents()
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C3"), "http://example.com")
$ws.Range("C3").ClearContents()
$ws.Hyperlinks.Delete()

# Restore selection/view state.
$ws.Range("M9").Select()
